$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Step 1: delete the old "Address" column (F) ---
$ws.Columns.Item(6).Delete() | Out-Null

# --- Step 2: insert "Parent / Guardian Name" before "Class & Section" (D) ---
$ws.Columns.Item(4).Insert() | Out-Null
$ws.Cells.Item(1, 4).Value = "Parent / Guardian Name"

# --- Step 3: insert "Date Of Birth" before "Blood Group" (F) ---
$ws.Columns.Item(6).Insert() | Out-Null
$ws.Cells.Item(1, 6).Value = "Date Of Birth"

# --- Step 4: append the new trailing columns (K..P) ---
$ws.Cells.Item(1, 11).Value = "Address Line - 1"
$ws.Cells.Item(1, 12).Value = "Address Line - 2"
$ws.Cells.Item(1, 13).Value = "Address Line - 3"
$ws.Cells.Item(1, 14).Value = "City"
$ws.Cells.Item(1, 15).Value = "Pincode"
$ws.Cells.Item(1, 16).Value = "Status"

# --- Step 5: apply header style to the new cells (bold/centered/filled, like the rest) ---
$headerStyleSource = $ws.Cells.Item(1, 1)
foreach ($col in 4,6,11,12,13,14,15,16) {
    $headerStyleSource.Copy() | Out-Null
    $ws.Cells.Item(1, $col).PasteSpecial(-4122) | Out-Null
}

# --- Step 6: column widths to match the final template layout ---
$ws.Columns.Item(1).ColumnWidth = 9
$ws.Columns.Item(2).ColumnWidth = 22
$ws.Columns.Item(3).ColumnWidth = 22
$ws.Columns.Item(4).ColumnWidth = 27.67
$ws.Columns.Item(5).ColumnWidth = 15.83
$ws.Columns.Item(6).ColumnWidth = 24.67
$ws.Columns.Item(7).ColumnWidth = 23.5
$ws.Columns.Item(8).ColumnWidth = 24.67
$ws.Columns.Item(9).ColumnWidth = 20.67
$ws.Columns.Item(10).ColumnWidth = 24.67
$ws.Columns.Item(11).ColumnWidth = 26.83
$ws.Columns.Item(12).ColumnWidth = 21.33
$ws.Columns.Item(13).ColumnWidth = 25.33
$ws.Columns.Item(14).ColumnWidth = 28.17
$ws.Columns.Item(15).ColumnWidth = 28.33
$ws.Columns.Item(16).ColumnWidth = 10.67
$ws.Columns.Item(17).ColumnWidth = 19.17

# --- Step 7: selection, matching the recorded cursor position in the template ---
$ws.Range("H18").Select() | Out-Null
